# The exercise's first three data rows (A1:B3 -> 159/36, 3359/35, 9072/8996) were
# removed, which shifts every row below them up by three. The former last row
# (A11:B11 -> 3577/91) ends up duplicated: once as the new last original row
# (new row 8) and once more as a freshly added row 9 with the same values.
#
# Net effect vs. the original 11 rows of data:
#   new row 1..8  == old rows 4..11
#   new row 9     == old row 11 (repeated)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 1-3 entirely; rows 4-11 shift up to become rows 1-8.
$ws.Range("A1:B3").EntireRow.Delete()

# Duplicate the (new) last row of data into a new row 9, matching the target.
$ws.Range("A9").Value2 = $ws.Range("A8").Value2
$ws.Range("B9").Value2 = $ws.Range("B8").Value2

# Mirror the saved selection state (column C, rows 1-9) from the target file.
$ws.Range("C1:C9").Select()
